# Update column F ("dSF") values for specific rows to reflect the
# repulled / recalculated data (commit: "repull data, push all data, mean calculation").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -3
    6  = -1
    7  = 3
    9  = 1
    13 = 3
    17 = 1
    18 = -1
    22 = 1
    28 = -1
    30 = 5
    34 = -2
    35 = -3
    41 = -8
    42 = -5
    46 = -1
    50 = 0
    51 = 1
    55 = 1
    56 = -3
    59 = 3
    61 = -7
    67 = -2
    70 = -5
    78 = 0
    79 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
